$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "PlantId"
$ws.Range("B1").Value = "ProductId"
$ws.Range("C1").Value = "Target"
$ws.Range("D1").Value = "Unit"

# --- Data rows (2-6) -----------------------------------------------------
$ws.Range("A2").Value = 71
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 42
$ws.Range("D2").Value = 2

$ws.Range("A3").Value = 77
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 41
$ws.Range("D3").Value = 2

$ws.Range("A4").Value = 78
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 2

$ws.Range("A5").Value = 79
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 39
$ws.Range("D5").Value = 2

$ws.Range("A6").Value = 75
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 2

# --- Remove the old 7th row entirely (was Plant_Id=11 / y / 60 / 1) -----
$ws.Rows.Item(7).Delete()

# --- Strip the leftover "applyAlignment=left" cell styling from A2:C6 ---
$ws.Range("A1:D6").ClearFormats()

# --- Selection, matching the saved sheet view ---------------------------
$ws.Range("C9").Select() | Out-Null
